# Updated symbol list on Sun Jan 29 06:12:46 UTC 2023 with GitHub Actions
#
# Refreshes the crypto ranking table (columns B..G, rows 2..51):
#  - Coin name / link shift up one rank for several mid-table rows
#  - Price (D) and Volume(1h) (E) values refreshed
#  - Hora (G) bumped from "5" to "6" for every data row
#
# D/E/G hold numeric-looking text (e.g. "0.0001301", "-1.27%", "6") that
# must stay literal strings (leading apostrophe forces text entry so
# Excel doesn't re-interpret/round them as numbers/percentages), then
# Style is reset to 'Normal' so no stray NumberFormat gets stamped onto
# the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = "'308.64"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'-1.27%"
$ws.Range('E2').Style = 'Normal'
$ws.Range('G2').Value = "'6"
$ws.Range('G2').Style = 'Normal'

# Row 3
$ws.Range('D3').Value = "'39.75"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'2.59%"
$ws.Range('E3').Style = 'Normal'
$ws.Range('G3').Value = "'6"
$ws.Range('G3').Style = 'Normal'

# Row 4
$ws.Range('D4').Value = "'5.129"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'0.36%"
$ws.Range('E4').Style = 'Normal'
$ws.Range('G4').Value = "'6"
$ws.Range('G4').Style = 'Normal'

# Row 5
$ws.Range('D5').Value = "'0.08135"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'-0.55%"
$ws.Range('E5').Style = 'Normal'
$ws.Range('G5').Value = "'6"
$ws.Range('G5').Style = 'Normal'

# Row 6
$ws.Range('D6').Value = "'1.948"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'-3.19%"
$ws.Range('E6').Style = 'Normal'
$ws.Range('G6').Value = "'6"
$ws.Range('G6').Style = 'Normal'

# Row 7
$ws.Range('D7').Value = "'8.151"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'2.97%"
$ws.Range('E7').Style = 'Normal'
$ws.Range('G7').Value = "'6"
$ws.Range('G7').Style = 'Normal'

# Row 8
$ws.Range('B8').Value = 'GateToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D8').Value = "'4.231"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'1.11%"
$ws.Range('E8').Style = 'Normal'
$ws.Range('G8').Value = "'6"
$ws.Range('G8').Style = 'Normal'

# Row 9
$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D9').Value = "'0.9289"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'-0.28%"
$ws.Range('E9').Style = 'Normal'
$ws.Range('G9').Value = "'6"
$ws.Range('G9').Style = 'Normal'

# Row 10
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D10').Value = "'0.1437"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'2.01%"
$ws.Range('E10').Style = 'Normal'
$ws.Range('G10').Value = "'6"
$ws.Range('G10').Style = 'Normal'

# Row 11
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').Value = "'0.1924"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'-1.17%"
$ws.Range('E11').Style = 'Normal'
$ws.Range('G11').Value = "'6"
$ws.Range('G11').Style = 'Normal'

# Row 12
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').Value = "'0.09129"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'-1.92%"
$ws.Range('E12').Style = 'Normal'
$ws.Range('G12').Value = "'6"
$ws.Range('G12').Style = 'Normal'

# Row 13
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').Value = "'0.03519"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'2.40%"
$ws.Range('E13').Style = 'Normal'
$ws.Range('G13').Value = "'6"
$ws.Range('G13').Style = 'Normal'

# Row 14
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').Value = "'0.09799"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'-0.51%"
$ws.Range('E14').Style = 'Normal'
$ws.Range('G14').Value = "'6"
$ws.Range('G14').Style = 'Normal'

# Row 15
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').Value = "'0.001394"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'-1.01%"
$ws.Range('E15').Style = 'Normal'
$ws.Range('G15').Value = "'6"
$ws.Range('G15').Style = 'Normal'

# Row 16
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').Value = "'0.005899"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'-1.52%"
$ws.Range('E16').Style = 'Normal'
$ws.Range('G16').Value = "'6"
$ws.Range('G16').Style = 'Normal'

# Row 17
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').Value = "'3.920"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'7.37%"
$ws.Range('E17').Style = 'Normal'
$ws.Range('G17').Value = "'6"
$ws.Range('G17').Style = 'Normal'

# Row 18
$ws.Range('D18').Value = "'3.366"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'-2.42%"
$ws.Range('E18').Style = 'Normal'
$ws.Range('G18').Value = "'6"
$ws.Range('G18').Style = 'Normal'

# Row 19
$ws.Range('E19').Value = "'-0.61%"
$ws.Range('E19').Style = 'Normal'
$ws.Range('G19').Value = "'6"
$ws.Range('G19').Style = 'Normal'

# Row 20
$ws.Range('D20').Value = "'0.1311"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'-0.13%"
$ws.Range('E20').Style = 'Normal'
$ws.Range('G20').Value = "'6"
$ws.Range('G20').Style = 'Normal'

# Row 21
$ws.Range('D21').Value = "'4.715"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'-1.77%"
$ws.Range('E21').Style = 'Normal'
$ws.Range('G21').Value = "'6"
$ws.Range('G21').Style = 'Normal'

# Row 22
$ws.Range('E22').Value = "'-1.15%"
$ws.Range('E22').Style = 'Normal'
$ws.Range('G22').Value = "'6"
$ws.Range('G22').Style = 'Normal'

# Row 23
$ws.Range('D23').Value = "'0.04378"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'-1.93%"
$ws.Range('E23').Style = 'Normal'
$ws.Range('G23').Value = "'6"
$ws.Range('G23').Style = 'Normal'

# Row 24
$ws.Range('E24').Value = "'-0.85%"
$ws.Range('E24').Style = 'Normal'
$ws.Range('G24').Value = "'6"
$ws.Range('G24').Style = 'Normal'

# Row 25
$ws.Range('D25').Value = "'0.004373"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'4.77%"
$ws.Range('E25').Style = 'Normal'
$ws.Range('G25').Value = "'6"
$ws.Range('G25').Style = 'Normal'

# Row 26
$ws.Range('D26').Value = "'0.0001301"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'-0.05%"
$ws.Range('E26').Style = 'Normal'
$ws.Range('G26').Value = "'6"
$ws.Range('G26').Style = 'Normal'

# Row 27
$ws.Range('D27').Value = "'0.0004000"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'-10.06%"
$ws.Range('E27').Style = 'Normal'
$ws.Range('G27').Value = "'6"
$ws.Range('G27').Style = 'Normal'

# Row 28
$ws.Range('G28').Value = "'6"
$ws.Range('G28').Style = 'Normal'

# Row 29
$ws.Range('G29').Value = "'6"
$ws.Range('G29').Style = 'Normal'

# Row 30
$ws.Range('G30').Value = "'6"
$ws.Range('G30').Style = 'Normal'

# Row 31
$ws.Range('G31').Value = "'6"
$ws.Range('G31').Style = 'Normal'

# Row 32
$ws.Range('G32').Value = "'6"
$ws.Range('G32').Style = 'Normal'

# Row 33
$ws.Range('G33').Value = "'6"
$ws.Range('G33').Style = 'Normal'

# Row 34
$ws.Range('G34').Value = "'6"
$ws.Range('G34').Style = 'Normal'

# Row 35
$ws.Range('G35').Value = "'6"
$ws.Range('G35').Style = 'Normal'

# Row 36
$ws.Range('G36').Value = "'6"
$ws.Range('G36').Style = 'Normal'

# Row 37
$ws.Range('G37').Value = "'6"
$ws.Range('G37').Style = 'Normal'

# Row 38
$ws.Range('G38').Value = "'6"
$ws.Range('G38').Style = 'Normal'

# Row 39
$ws.Range('D39').Value = "'0.02057"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'-3.83%"
$ws.Range('E39').Style = 'Normal'
$ws.Range('G39').Value = "'6"
$ws.Range('G39').Style = 'Normal'

# Row 40
$ws.Range('D40').Value = "'0.05051"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'-2.61%"
$ws.Range('E40').Style = 'Normal'
$ws.Range('G40').Value = "'6"
$ws.Range('G40').Style = 'Normal'

# Row 41
$ws.Range('D41').Value = "'0.007437"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'-0.28%"
$ws.Range('E41').Style = 'Normal'
$ws.Range('G41').Value = "'6"
$ws.Range('G41').Style = 'Normal'

# Row 42
$ws.Range('D42').Value = "'0.009865"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'-1.39%"
$ws.Range('E42').Style = 'Normal'
$ws.Range('G42').Value = "'6"
$ws.Range('G42').Style = 'Normal'

# Row 43
$ws.Range('E43').Value = "'-0.32%"
$ws.Range('E43').Style = 'Normal'
$ws.Range('G43').Value = "'6"
$ws.Range('G43').Style = 'Normal'

# Row 44
$ws.Range('D44').Value = "'0.002131"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'-0.04%"
$ws.Range('E44').Style = 'Normal'
$ws.Range('G44').Value = "'6"
$ws.Range('G44').Style = 'Normal'

# Row 45
$ws.Range('D45').Value = "'0.008630"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'-11.96%"
$ws.Range('E45').Style = 'Normal'
$ws.Range('G45').Value = "'6"
$ws.Range('G45').Style = 'Normal'

# Row 46
$ws.Range('D46').Value = "'0.00006364"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'0.69%"
$ws.Range('E46').Style = 'Normal'
$ws.Range('G46').Value = "'6"
$ws.Range('G46').Style = 'Normal'

# Row 47
$ws.Range('E47').Value = "'-0.08%"
$ws.Range('E47').Style = 'Normal'
$ws.Range('G47').Value = "'6"
$ws.Range('G47').Style = 'Normal'

# Row 48
$ws.Range('D48').Value = "'0.002713"
$ws.Range('D48').Style = 'Normal'
$ws.Range('G48').Value = "'6"
$ws.Range('G48').Style = 'Normal'

# Row 49
$ws.Range('E49').Value = "'-18.81%"
$ws.Range('E49').Style = 'Normal'
$ws.Range('G49').Value = "'6"
$ws.Range('G49').Style = 'Normal'

# Row 50
$ws.Range('E50').Value = "'-0.08%"
$ws.Range('E50').Style = 'Normal'
$ws.Range('G50').Value = "'6"
$ws.Range('G50').Style = 'Normal'

# Row 51
$ws.Range('E51').Value = "'-0.08%"
$ws.Range('E51').Style = 'Normal'
$ws.Range('G51').Value = "'6"
$ws.Range('G51').Style = 'Normal'
